$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Card" sample table originally held just Value/State. Synoptic Panel's
# glitch was fixed by extending the sample data with the six Great/Moderate/Fail
# threshold columns used to drive the state coloring, so insert 6 new columns
# after State (this also shifts the stray column-C width formatting to column I,
# matching a real "insert columns" user action).
$ws.Range("C1:H1").EntireColumn.Insert()

# Grow the Table2 list object so the new columns become real table columns
# (not just bare cells) and the table/autofilter ranges stay in sync.
$lo = $ws.ListObjects(1)
$lo.Resize($ws.Range("A1:H2"))

# New header row text
$ws.Range("C1").Value = "Great from"
$ws.Range("D1").Value = "Great to"
$ws.Range("E1").Value = "Moderate from"
$ws.Range("F1").Value = "Moderate to"
$ws.Range("G1").Value = "Fail from"
$ws.Range("H1").Value = "Fail to"

# New data row values (State changes from 1 to 40; new threshold values added)
$ws.Range("B2").Value = 40
$ws.Range("C2").Value = 50
$ws.Range("D2").Value = 1048576
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 50
$ws.Range("G2").Value = -1048576
$ws.Range("H2").Value = 20

# Matches the saved selection left behind in the authored workbook
$ws.Range("K16").Select()
